$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the header row "municipio", "CASOS", "óbitos"); this shifts
# all subsequent data rows up by one.
$ws.Rows.Item(2).Delete()
